$d = $word.ActiveDocument

$replacements = @(
    @("3+37=", "28+60="),
    @("69+24=", "42+24="),
    @("81-75=", "50+29="),
    @("71-6=", "56+19="),
    @("52-10=", "69+10="),
    @("94-21=", "60-34="),
    @("77-58=", "10+12="),
    @("49-44=", "52+33="),
    @("59-35=", "45+35="),
    @("49+2=", "16-7="),
    @("41-9=", "50+44="),
    @("43-14=", "26+31="),
    @("33-28=", "42+30="),
    @("36+53=", "6+33="),
    @("75-70=", "54-38="),
    @("33+7=", "86+1="),
    @("10+68=", "95-88="),
    @("17+33=", "28+52="),
    @("88-68=", "73-51="),
    @("24+12=", "34+35="),
    @("41+29=", "76+5="),
    @("63-59=", "93+2="),
    @("78-57=", "43+16="),
    @("86+12=", "93-55="),
    @("69+3=", "33+53="),
    @("55+6=", "62-47="),
    @("11+67=", "87-65="),
    @("48+12=", "36+32="),
    @("10+78=", "48-23="),
    @("95-65=", "58-58="),
    @("73+12=", "38+14="),
    @("40+21=", "93-46="),
    @("82+13=", "86-6="),
    @("36+0=", "33+32="),
    @("28+18=", "85+6="),
    @("97-46=", "92-11="),
    @("1+4=", "24-11="),
    @("46-45=", "79-1="),
    @("0+13=", "76-54="),
    @("86-85=", "59-30="),
    @("89-1=", "54+15="),
    @("43-22=", "43+38="),
    @("25+64=", "93-81="),
    @("99-44=", "49-49="),
    @("48-30=", "51+5="),
    @("50-31=", "21+11="),
    @("2+79=", "61-24="),
    @("67+3=", "74+7="),
    @("98-94=", "30+4="),
    @("5-1=", "50-46="),
    @("92-40=", "98-2="),
    @("29+55=", "75-41="),
    @("87-24=", "16+70="),
    @("32+53=", "37-16="),
    @("22+59=", "39+6="),
    @("76+23=", "2+47="),
    @("12+58=", "82-27="),
    @("5+23=", "6+18="),
    @("87+0=", "64-48="),
    @("60-58=", "15+45="),
    @("5+11=", "55-33="),
    @("36-7=", "18+43="),
    @("14+42=", "98-24="),
    @("57+32=", "52+20="),
    @("97-15=", "85-7="),
    @("32+0=", "89+2="),
    @("63+34=", "98+1="),
    @("62+10=", "1+54="),
    @("58-8=", "54-13="),
    @("96-68=", "37+47="),
    @("75-25=", "54-42="),
    @("19+27=", "29-22="),
    @("81-17=", "35-6="),
    @("39-8=", "34-32="),
    @("81+1=", "24+55="),
    @("98-4=", "42+35="),
    @("93-92=", "24+27="),
    @("24+66=", "49-12="),
    @("86-0=", "85-58="),
    @("29+15=", "44+13="),
    @("17+39=", "42-21="),
    @("78-5=", "21+78="),
    @("27+59=", "86-8="),
    @("53-1=", "27-15="),
    @("42-39=", "56+2="),
    @("52+40=", "98-36="),
    @("25+65=", "44-36="),
    @("43-2=", "61+4="),
    @("71-60=", "33+59="),
    @("41+9=", "83+14="),
    @("92-23=", "66-7="),
    @("46+12=", "76-32="),
    @("11+76=", "45+9="),
    @("61-16=", "9+2="),
    @("70-53=", "34+39="),
    @("25+21=", "45+45="),
    @("70+22=", "30+60="),
    @("3+32=", "83-48="),
    @("43-8=", "93-18="),
    @("86-78=", "0+88="),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $result = $find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $result) {
        Write-Host "WARNING: replacement failed for $old -> $new"
    }
}

Write-Host "Done."